$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4360621853365316
$ws.Range("C2").Value = 0.6094159519752258
$ws.Range("D2").Value = 0.6952771648018571
$ws.Range("E2").Value = 0.7487255867942739
$ws.Range("B3").Value = 0.450039492624489
$ws.Range("C3").Value = 0.6205207970795291
$ws.Range("D3").Value = 0.7012035129671159
$ws.Range("E3").Value = 0.7524391359638646
$ws.Range("B4").Value = 0.4139473131987446
$ws.Range("C4").Value = 0.5928651984725631
$ws.Range("D4").Value = 0.6869053491050978
$ws.Range("E4").Value = 0.7427761591611421
$ws.Range("B5").Value = 0.473703034265991
$ws.Range("C5").Value = 0.6401494306782807
$ws.Range("D5").Value = 0.7044098395282327
$ws.Range("E5").Value = 0.7593728668952565
$ws.Range("B6").Value = 0.4689871014179974
$ws.Range("C6").Value = 0.6360887822946704
$ws.Range("D6").Value = 0.7024036911943025
$ws.Range("E6").Value = 0.7579578215545946
$ws.Range("B7").Value = 0.4807044994505518
$ws.Range("C7").Value = 0.6441383609149771
$ws.Range("D7").Value = 0.7127657133716443
$ws.Range("E7").Value = 0.7598586426095999
$ws.Range("B8").Value = 0.4166872841406259
$ws.Range("C8").Value = 0.5873098263957595
$ws.Range("D8").Value = 0.6845477254926853
$ws.Range("E8").Value = 0.7388605506954972
$ws.Range("B9").Value = 0.4753935026063243
$ws.Range("C9").Value = 0.6392664312515203
$ws.Range("D9").Value = 0.7101543221519784
$ws.Range("E9").Value = 0.7583282497459432
$ws.Range("B10").Value = 0.4848457968482907
$ws.Range("C10").Value = 0.6452139183127141
$ws.Range("D10").Value = 0.7122758823361749
$ws.Range("E10").Value = 0.758778634037358
$ws.Range("B11").Value = 0.4844162104220959
$ws.Range("C11").Value = 0.6454645891659274
$ws.Range("D11").Value = 0.712790944671245
$ws.Range("E11").Value = 0.7592666848883107
$ws.Range("B12").Value = 0.4749351916719513
$ws.Range("C12").Value = 0.6356261515511745
$ws.Range("D12").Value = 0.7060317019685548
$ws.Range("E12").Value = 0.7525739072257516
$ws.Range("B13").Value = 0.4838181537396462
$ws.Range("C13").Value = 0.6440407900987742
$ws.Range("D13").Value = 0.7113333033618732
$ws.Range("E13").Value = 0.7578752676462324
